# Saved progress at the end of the loop
# Updates the "Qty executed upto date" (C) and the derived "Upto date
# Amount" / "Amount Since prev bill" (G/H) columns for the bill-summary
# rows, then refreshes the Grand Total / Net Payable rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving it as TEXT (the G/H
# "Upto date Amount" cells are stored as text strings like "15872.00",
# not numbers). Going through a self-referencing text formula and then
# collapsing it to a static value (PasteSpecial -> values only) keeps the
# cell's text type without leaving a formula behind and without having to
# touch the cell's number format / style.
function Set-TextValue {
    param($addr, [string]$text)
    $r = $ws.Range($addr)
    $escaped = $text -replace '"', '""'
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# Qty executed upto date (plain numbers)
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 62
$ws.Range("C10").Value = 96
$ws.Range("C11").Value = 25
$ws.Range("C12").Value = 59
$ws.Range("C13").Value = 20
$ws.Range("C14").Value = 10
$ws.Range("C15").Value = 9
$ws.Range("C16").Value = 83
$ws.Range("C17").Value = 37

# Upto date Amount (text-formatted numbers)
Set-TextValue "G9"  "15872.00"
Set-TextValue "G10" "45312.00"
Set-TextValue "G11" "16550.00"
Set-TextValue "G13" "2720.00"
Set-TextValue "G14" "230.00"

# Grand Total / Net Payable rows
Set-TextValue "G19" "80684.00"
Set-TextValue "H19" "80684.00"
Set-TextValue "G21" "80684.00"
Set-TextValue "H21" "80684.00"
